$wb = $excel.ActiveWorkbook

# --- Sheet 1 "展览": update "想去人数" (F column) counts for several events ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 172
$ws1.Range("F3").Value = 188
$ws1.Range("F4").Value = 2933
$ws1.Range("F5").Value = 205
$ws1.Range("F6").Value = 106
$ws1.Range("F7").Value = 190
$ws1.Range("F8").Value = 1619
$ws1.Range("F9").Value = 1601
$ws1.Range("F10").Value = 51
$ws1.Range("F14").Value = 185
$ws1.Range("F15").Value = 22
$ws1.Range("F16").Value = 221
$ws1.Range("F18").Value = 224
$ws1.Range("F21").Value = 34
$ws1.Range("F22").Value = 3
$ws1.Range("F23").Value = 343
$ws1.Range("F24").Value = 121
$ws1.Range("F25").Value = 91
$ws1.Range("F26").Value = 13
$ws1.Range("F27").Value = 1941
$ws1.Range("F29").Value = 448
$ws1.Range("F30").Value = 8
$ws1.Range("F31").Value = 146
$ws1.Range("F32").Value = 572
$ws1.Range("F34").Value = 325
$ws1.Range("F36").Value = 478

# --- Sheet 2 "演出": the single listed event ("南昌·六一Big Day...") is gone ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()

# --- Sheet 4 "全部类型": drop the duplicated "南昌·六一Big Day" rows, add the missing
#     "宜春·COMIC WORLD..." row and refresh the same F-column counts as sheet 1 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Delete()
$ws4.Rows.Item(2).Delete()

$ws4.Rows.Item(22).Insert()
$ws4.Range("A21:I21").Copy()
$ws4.Range("A22:I22").PasteSpecial(-4122)
$ws4.Range("B22").NumberFormat = "@"
$ws4.Range("A22").Value = 21
$ws4.Range("B22").Value = "2024-07-13"
$ws4.Range("C22").Value = "宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华"
$ws4.Range("D22").Value = "宜春国际商贸城会展中心 宜春国际商贸城会展中心"
$ws4.Range("E22").Value = "2024.07.13 10:00-07.14 17:00"
$ws4.Range("F22").Value = 3
$ws4.Range("G22").Value = 55
$ws4.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=86667"
$ws4.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg"

$ws4.Range("F2").Value = 172
$ws4.Range("F3").Value = 188
$ws4.Range("F4").Value = 2933
$ws4.Range("F5").Value = 205
$ws4.Range("F6").Value = 106
$ws4.Range("F7").Value = 190
$ws4.Range("F8").Value = 1619
$ws4.Range("F9").Value = 1601
$ws4.Range("F10").Value = 51
$ws4.Range("F14").Value = 185
$ws4.Range("F15").Value = 22
$ws4.Range("F16").Value = 221
$ws4.Range("F18").Value = 224
$ws4.Range("F21").Value = 34
$ws4.Range("F23").Value = 343
$ws4.Range("F24").Value = 121
$ws4.Range("F25").Value = 91
$ws4.Range("F26").Value = 13
$ws4.Range("F27").Value = 1941
$ws4.Range("F29").Value = 448
$ws4.Range("F30").Value = 8
$ws4.Range("F31").Value = 146
$ws4.Range("F32").Value = 572
$ws4.Range("F34").Value = 325
$ws4.Range("F36").Value = 478
